# Generate Report for Handoff
# -----------------------------------------------------------------------
# The localization-status report moved from "In Translation" to
# "Ready for handoff" and the handoff timestamps were refreshed. Update
# the Overview roll-up sheet plus the two per-locale detail sheets
# (zh-cn, de-de) to match, and widen the Status-ish date columns that
# Excel re-flowed to fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$ws1 = $wb.Sheets.Item("Overview")
$ws1.Range("E2").Value = "Ready for handoff"
$ws1.Range("F2").Value = "Ready for handoff"
$ws1.Range("G2").Value = "2016-09-05 13:12:58"

# Widen the zh-cn / de-de status columns on the overview sheet (Excel
# quantizes column widths to whole pixels, so this lands on the closest
# representable width to the authored 17.2159881591797).
$ws1.Columns.Item(5).ColumnWidth = 16.3
$ws1.Columns.Item(6).ColumnWidth = 16.3

# --- zh-cn detail sheet ------------------------------------------------
$ws2 = $wb.Sheets.Item("zh-cn")
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("H2").Value = "2016-09-05 13:12:53"
$ws2.Columns.Item(3).ColumnWidth = 16.3

# --- de-de detail sheet ------------------------------------------------
$ws3 = $wb.Sheets.Item("de-de")
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("H2").Value = "2016-09-05 13:12:58"
$ws3.Columns.Item(3).ColumnWidth = 16.3
